# edit.ps1 - applies the two changes captured by the commit's OOXML diff:
#
#   1. The table on slide 5 gets its table style switched from the
#      custom "Table_0" style ({069CDA57-0D43-47A6-8A42-7AA7F0E96C46})
#      to the built-in "No Style, Table Grid" style
#      ({ECB931A1-7BFD-4CCE-BB0F-C497967DC934}).
#
#   2. The presentation's theme (ppt/theme/theme1.xml, the one actually
#      used by the slide master / slides) has its 12 scheme colors
#      changed from the "Integral" (Red Violet) palette over to the
#      stock "Office Theme" palette. (Font scheme / format scheme were
#      already identical between the two themes in this deck, so only
#      the color scheme actually needs to move.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table on slide 5.
# ---------------------------------------------------------------------
$targetStyleId = "{ECB931A1-7BFD-4CCE-BB0F-C497967DC934}"
$tableFound = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            # Table styles can't be assigned through the .Style property -
            # PowerPoint requires ApplyStyle(styleId) instead.
            $shape.Table.ApplyStyle($targetStyleId)
            $tableFound = $true
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme color scheme over to the stock "Office Theme"
#    palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in that
#    order, matching ThemeColorScheme.Colors(1..12)).
# ---------------------------------------------------------------------
$officeThemeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    # PowerPoint's RGB storage is packed as 0xBBGGRR.
    $packed = ($b * 0x10000) + ($g * 0x100) + $r
    $colorScheme.Colors($i).RGB = $packed
}

Write-Host "tableFound=$tableFound"
